# "NCC with feature data error fix on the eu find distance"
# The NCC-with-feature Euclidean distance values (row 1, columns B:G) were
# recomputed; update the cells accordingly (this also drops the old centered
# style that was applied to that row, matching the corrected data range).
# The chart was also repositioned from the right-hand side of the sheet to
# sit directly under the data table, and the active selection moved on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (NCC with Feature): corrected Euclidean-distance values ---
$ws.Range("B1:G1").Style = "Normal"
$ws.Range("B1").Value = 0.18909089850313116
$ws.Range("C1").Value = 0.1883158691412665
$ws.Range("D1").Value = 0.18950428968160471
$ws.Range("E1").Value = 0.18206111442472653
$ws.Range("F1").Value = 0.18689628297044972
$ws.Range("G1").Value = 0.18708011148677436

# --- Reposition the chart: from col I/row 4-ish to col A/row 8-ish, under the table ---
$co = $ws.ChartObjects(1)
$co.Left = 2.25
$co.Top = 118.12488188976378
$co.Width = 541.9112115526575
$co.Height = 295.65

# --- Move the active selection ---
[void]$ws.Range("K19").Select()
